$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 25.5
$ws.Range("C3").Value = 208
$ws.Range("D3").Value = 224
$ws.Range("C4").Value = 154
$ws.Range("D4").Value = 166
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 33
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 63
$ws.Range("C7").Value = 85
$ws.Range("D7").Value = 95
$ws.Range("C8").Value = 111
$ws.Range("D8").Value = 120.5
$ws.Range("C10").Value = 256
$ws.Range("D10").Value = 269
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 74
$ws.Range("C12").Value = 54
$ws.Range("D12").Value = 69
$ws.Range("C13").Value = 17
$ws.Range("D13").Value = 30
$ws.Range("C14").Value = 125
$ws.Range("D14").Value = 149
$ws.Range("C15").Value = 105
$ws.Range("D15").Value = 112.5
$ws.Range("C17").Value = 185
$ws.Range("D17").Value = 193
$ws.Range("C19").Value = 100
$ws.Range("D19").Value = 97.5
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 37
$ws.Range("C21").Value = 114
$ws.Range("D21").Value = 114.5
$ws.Range("C22").Value = 864
$ws.Range("D22").Value = 864.5
$ws.Range("C23").Value = 231
$ws.Range("D23").Value = 199
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 33.5
$ws.Range("C25").Value = 174
$ws.Range("D25").Value = 186.5
$ws.Range("C26").Value = 210
$ws.Range("D26").Value = 227
$ws.Range("C27").Value = 102
$ws.Range("D27").Value = 102
$ws.Range("C28").Value = 160
$ws.Range("D28").Value = 191.5
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 66
$ws.Range("C30").Value = 38
$ws.Range("D30").Value = 37.5
$ws.Range("C31").Value = 49
$ws.Range("D31").Value = 53
$ws.Range("C32").Value = 161
$ws.Range("D32").Value = 191.5
$ws.Range("C33").Value = 77
$ws.Range("D33").Value = 82.5
$ws.Range("C34").Value = 94
$ws.Range("D34").Value = 111.5
$ws.Range("C35").Value = 75
$ws.Range("D35").Value = 77.5
$ws.Range("C36").Value = 43
$ws.Range("D36").Value = 48.5
$ws.Range("C37").Value = 39
$ws.Range("D37").Value = 46
$ws.Range("C38").Value = 45
$ws.Range("D38").Value = 59.5
$ws.Range("C39").Value = 10
$ws.Range("D39").Value = 31
$ws.Range("C40").Value = 922
$ws.Range("D40").Value = 916
$ws.Range("C41").Value = 71
$ws.Range("D41").Value = 80
$ws.Range("C42").Value = 36
$ws.Range("D42").Value = 37
$ws.Range("C43").Value = 59
$ws.Range("D43").Value = 66
$ws.Range("C44").Value = 334
$ws.Range("D44").Value = 345
$ws.Range("C45").Value = 83
$ws.Range("D45").Value = 85
$ws.Range("C46").Value = 30
$ws.Range("D46").Value = 31.5
$ws.Range("C47").Value = 280
$ws.Range("D47").Value = 274
$ws.Range("C48").Value = 212
$ws.Range("D48").Value = 212.5
$ws.Range("C50").Value = 22
$ws.Range("D50").Value = 31.5
$ws.Range("C51").Value = 336
$ws.Range("D51").Value = 347
$ws.Range("C52").Value = 25
$ws.Range("D52").Value = 24.5
$ws.Range("C54").Value = 35
$ws.Range("D54").Value = 40.5
$ws.Range("C55").Value = 189
$ws.Range("D55").Value = 182
$ws.Range("C56").Value = 118
$ws.Range("D56").Value = 125
$ws.Range("C57").Value = 188
$ws.Range("D57").Value = 195
$ws.Range("C59").Value = 9
$ws.Range("D59").Value = 66
$ws.Range("C60").Value = 88
$ws.Range("D60").Value = 97
$ws.Range("C61").Value = 63
$ws.Range("D61").Value = 66.5
$ws.Range("C62").Value = 371
$ws.Range("D62").Value = 371
$ws.Range("C63").Value = 61
$ws.Range("D63").Value = 66.5
$ws.Range("C64").Value = 270
$ws.Range("D64").Value = 265.5
$ws.Range("C65").Value = 169
$ws.Range("D65").Value = 179.5
$ws.Range("C66").Value = 12
$ws.Range("D66").Value = 69.5
$ws.Range("C67").Value = 20
$ws.Range("D67").Value = 25.5
$ws.Range("C68").Value = 57
$ws.Range("D68").Value = 61.5
$ws.Range("C69").Value = 123
$ws.Range("D69").Value = 128.5
$ws.Range("C70").Value = 65
$ws.Range("D70").Value = 71
$ws.Range("C71").Value = 29
$ws.Range("D71").Value = 30.5
$ws.Range("C72").Value = 116
$ws.Range("D72").Value = 128
$ws.Range("C73").Value = 51
$ws.Range("D73").Value = 55
$ws.Range("C75").Value = 69
$ws.Range("D75").Value = 72.5
$ws.Range("C76").Value = 125.6081081081081
